# The deck originally ships with the "Integral" theme (ppt/theme/theme1.xml,
# used by the slide master) and an "Office Theme" (ppt/theme/theme2.xml,
# used by the notes master). The authored change swaps the two themes'
# content: the slide-master-facing theme becomes the stock "Office Theme"
# palette (and, symmetrically, the notes-master-facing theme becomes the
# "Integral" palette). This script applies the reachable half of that swap
# through the PowerPoint object model: it recolors the active design's
# theme color scheme (ppt/theme/theme1.xml) from "Integral" to the twelve
# standard "Office Theme" colors.

$p = $ppt.ActivePresentation
$sm = $p.SlideMaster
$themeColors = $sm.Theme.ThemeColorScheme

function ToBGR([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Index order per MsoThemeColorSchemeIndex: dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink.
$officeThemeHex = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = ToBGR($officeThemeHex[$i - 1])
}
